# Update "想去人数" (want-to-go count) figures for the 南宁-漫展信息 workbook.
# The same underlying data is duplicated on the "展览" sheet and the
# "全部类型" sheet, so both need to be updated identically.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 541
    $ws.Range("F3").Value = 3499
    $ws.Range("F5").Value = 682
}
